$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '56.475.84'
$ws.Range('E2').Value = '  -2.04%  '

# Row 3
$ws.Range('D3').Value = '2.995.62'
$ws.Range('E3').Value = '  -4.18%  '

# Row 4
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
$ws.Range('D5').Value = '495.74'
$ws.Range('E5').Value = '  -5.08%  '

# Row 6
$ws.Range('D6').Value = '134.36'
$ws.Range('E6').Value = '  -0.33%  '

# Row 7
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.11%  '

# Row 8
$ws.Range('D8').Value = '2.990.44'
$ws.Range('E8').Value = '  -4.27%  '

# Row 9
$ws.Range('D9').Value = '0.425'
$ws.Range('E9').Value = '  -4.27%  '

# Row 10
$ws.Range('D10').Value = '7.20'
$ws.Range('E10').Value = '  -0.22%  '

# Row 11
$ws.Range('D11').Value = '0.104'
$ws.Range('E11').Value = '  -5.42%  '

# Row 12
$ws.Range('D12').Value = '0.353'
$ws.Range('E12').Value = '  -8.10%  '

# Row 13
$ws.Range('E13').Value = '  +0.57%  '

# Row 14
$ws.Range('D14').Value = '3.503.42'
$ws.Range('E14').Value = '  -4.07%  '

# Row 15
$ws.Range('D15').Value = '25.01'
$ws.Range('E15').Value = '  -1.82%  '

# Row 16
$ws.Range('D16').Value = '56.314.91'
$ws.Range('E16').Value = '  -2.18%  '

# Row 17
$ws.Range('D17').Value = '2.995.71'
$ws.Range('E17').Value = '  -3.79%  '

# Row 18
$ws.Range('D18').Value = '0.0000144'
$ws.Range('E18').Value = '  -4.91%  '

# Row 19
$ws.Range('D19').Value = '5.84'
$ws.Range('E19').Value = '  +0.86%  '

# Row 20
$ws.Range('D20').Value = '12.38'
$ws.Range('E20').Value = '  -4.85%  '

# Row 21
$ws.Range('D21').Value = '7.71'
$ws.Range('E21').Value = '  -3.74%  '

# Row 22
$ws.Range('D22').Value = '322.96'
$ws.Range('E22').Value = '  -6.05%  '

# Row 23
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.10%  '

# Row 24
$ws.Range('D24').Value = '0.464'
$ws.Range('E24').Value = '  -7.96%  '

# Row 25
$ws.Range('D25').Value = '61.09'
$ws.Range('E25').Value = '  -10.58%  '

# Row 26
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  -0.22%  '

# Row 27
$ws.Range('E27').Value = '  -1.28%  '

# Row 28
$ws.Range('D28').Value = '0.0₃0884'
$ws.Range('E28').Value = '  -6.57%  '

# Row 29
$ws.Range('E29').Value = '  +0.12%  '

# Row 30
$ws.Range('D30').Value = '6.58'
$ws.Range('E30').Value = '  -2.91%  '

# Row 31
$ws.Range('D31').Value = '6.73'
$ws.Range('E31').Value = '  -2.15%  '

# Row 32
$ws.Range('D32').Value = '1.18'

# Row 33
$ws.Range('E33').Value = '  -7.86%  '

# Row 34
$ws.Range('D34').Value = '19.86'
$ws.Range('E34').Value = '  -7.75%  '

# Row 35
$ws.Range('D35').Value = '151.17'
$ws.Range('E35').Value = '  -4.38%  '

# Row 36
$ws.Range('D36').Value = '4.50'
$ws.Range('E36').Value = '  -5.71%  '

# Row 37
$ws.Range('E37').Value = '  -6.69%  '

# Row 38
$ws.Range('D38').Value = '5.65'
$ws.Range('E38').Value = '  -8.43%  '

# Row 39
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.0662'
$ws.Range('E39').Value = '  -4.11%  '

# Row 40
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').Value = '23.59'
$ws.Range('E40').Value = '  -6.60%  '

# Row 41
$ws.Range('D41').Value = '3.025.91'
$ws.Range('E41').Value = '  -3.90%  '

# Row 42
$ws.Range('D42').Value = '36.81'
$ws.Range('E42').Value = '  -8.43%  '

# Row 43
$ws.Range('E43').Value = '  +0.23%  '

# Row 44
$ws.Range('E44').Value = '  -4.38%  '

# Row 45
$ws.Range('D45').Value = '0.640'
$ws.Range('E45').Value = '  -6.16%  '

# Row 46
$ws.Range('D46').Value = '1.42'
$ws.Range('E46').Value = '  -2.23%  '

# Row 47
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').Value = '3.56'
$ws.Range('E47').Value = '  -8.71%  '

# Row 48
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.172.47'
$ws.Range('E48').Value = '  -3.49%  '

# Row 49
$ws.Range('D49').Value = '0.0238'
$ws.Range('E49').Value = '  +2.04%  '

# Row 50
$ws.Range('D50').Value = '19.38'
$ws.Range('E50').Value = '  -3.44%  '

# Row 51
$ws.Range('D51').Value = '1.85'
$ws.Range('E51').Value = '  +2.00%  '
